$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns AZ (52), BA (53), BB (54)
$ws.Range("AZ1").Value = "m_adult_literacy_pct"
$ws.Range("BA1").Value = "m_homicides_per_100k"
$ws.Range("BB1").Value = "m_tax_revenue_pct_gdp"

# Copy the header style from an existing header cell (AY1) to the new header cells
$ws.Range("AY1").Copy()
$ws.Range("AZ1:BB1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$values = @{
    2  = @(0, 0, 1)
    3  = @(1, 1, 1)
    4  = @(1, 1, 1)
    5  = @(1, 1, 0)
    6  = @(1, 1, 0)
    7  = @(1, 0, 1)
    8  = @(1, 1, 0)
    9  = @(1, 0, 1)
    10 = @(1, 1, 1)
    11 = @(0, 1, 0)
    12 = @(1, 1, 1)
    13 = @(0, 0, 0)
    14 = @(1, 1, 1)
    15 = @(1, 1, 1)
}

foreach ($row in $values.Keys) {
    $vals = $values[$row]
    $ws.Range("AZ$row").Value = $vals[0]
    $ws.Range("BA$row").Value = $vals[1]
    $ws.Range("BB$row").Value = $vals[2]
}
